$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.262.88"
$ws.Range("E2").Value = "'  -2.06%  "
$ws.Range("D3").Value = "'2.582.28"
$ws.Range("E3").Value = "'  -2.14%  "
$ws.Range("E4").Value = "'  -0.06%  "
$ws.Range("D5").Value = "'562.76"
$ws.Range("E5").Value = "'  -1.41%  "
$ws.Range("D6").Value = "'142.64"
$ws.Range("E6").Value = "'  -2.81%  "
$ws.Range("E7").Value = "'  +0.29%  "
$ws.Range("E8").Value = "'  -2.28%  "
$ws.Range("D9").Value = "'2.590.36"
$ws.Range("E9").Value = "'  -2.74%  "
$ws.Range("D10").Value = "'6.64"
$ws.Range("E10").Value = "'  -2.70%  "
$ws.Range("D11").Value = "'0.104"
$ws.Range("E11").Value = "'  -0.69%  "
$ws.Range("E12").Value = "'  +11.64%  "
$ws.Range("D13").Value = "'0.355"
$ws.Range("E13").Value = "'  +3.52%  "
$ws.Range("D14").Value = "'3.035.97"
$ws.Range("E14").Value = "'  -2.36%  "
$ws.Range("D15").Value = "'23.25"
$ws.Range("E15").Value = "'  +7.00%  "
$ws.Range("D16").Value = "'59.222.23"
$ws.Range("E16").Value = "'  -2.03%  "
$ws.Range("E17").Value = "'  -0.35%  "
$ws.Range("D18").Value = "'2.593.29"
$ws.Range("E18").Value = "'  -2.25%  "
$ws.Range("E19").Value = "'  +0.47%  "
$ws.Range("D20").Value = "'337.20"
$ws.Range("E20").Value = "'  -2.33%  "
$ws.Range("D21").Value = "'10.36"
$ws.Range("E21").Value = "'  -0.65%  "
$ws.Range("E22").Value = "'  -0.06%  "
$ws.Range("E23").Value = "'  +0.19%  "
$ws.Range("D24").Value = "'64.17"
$ws.Range("E24").Value = "'  -4.05%  "
$ws.Range("D25").Value = "'0.467"
$ws.Range("E25").Value = "'  +5.42%  "
$ws.Range("E26").Value = "'  +0.50%  "
$ws.Range("E27").Value = "'  -3.03%  "
$ws.Range("E28").Value = "'  -0.54%  "
$ws.Range("D29").Value = "'0.0₃0774"
$ws.Range("E29").Value = "'  -0.45%  "
$ws.Range("E31").Value = "'  -2.88%  "
$ws.Range("D32").Value = "'6.11"
$ws.Range("E32").Value = "'  -0.03%  "
$ws.Range("D33").Value = "'159.86"
$ws.Range("E33").Value = "'  +2.35%  "
$ws.Range("D34").Value = "'18.97"
$ws.Range("E34").Value = "'  -1.30%  "
$ws.Range("E35").Value = "'  -1.61%  "
$ws.Range("D36").Value = "'1.17"
$ws.Range("E36").Value = "'  -1.01%  "
$ws.Range("E37").Value = "'  -3.61%  "
$ws.Range("D38").Value = "'0.872"
$ws.Range("E38").Value = "'  -4.50%  "
$ws.Range("D39").Value = "'37.51"
$ws.Range("E39").Value = "'  -0.38%  "
$ws.Range("E40").Value = "'  -2.37%  "
$ws.Range("E41").Value = "'  +0.08%  "
$ws.Range("D42").Value = "'293.64"
$ws.Range("E42").Value = "'  -4.55%  "
$ws.Range("B43").Value = "'Aave"
$ws.Range("C43").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'132.21"
$ws.Range("E43").Value = "'  +5.49%  "
$ws.Range("B44").Value = "'FirstDigitalUSD"
$ws.Range("C44").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "'  +0.50%  "
$ws.Range("E45").Value = "'  -0.93%  "
$ws.Range("D46").Value = "'0.595"
$ws.Range("E46").Value = "'  -2.09%  "
$ws.Range("E47").Value = "'  -0.13%  "
$ws.Range("E48").Value = "'  -2.66%  "
$ws.Range("D49").Value = "'19.01"
$ws.Range("E49").Value = "'  -2.30%  "
$ws.Range("E50").Value = "'  -1.00%  "
$ws.Range("D51").Value = "'18.67"
$ws.Range("E51").Value = "'  -0.05%  "
